# Adds a new "Player Info" sheet ahead of "ODI Batting", and reworks the
# MATCH_CARD_LINK columns (full scorecard URL) into MATCH_CODE columns
# (just the numeric match code) on "ODI Batting" / "ODI Bowling".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet before "ODI Batting" so the tab
#    order becomes: Player Info, ODI Batting, ODI Bowling.
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item("ODI Batting"))
$playerInfo.Name = "Player Info"

# Re-fetch "ODI Batting" by name - the sheet reference used above to
# position the new sheet becomes stale (it now points at the newly
# inserted sheet instead) once the insertion shifts indices.
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Header row.
$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

# Reuse the same header look (bold, centered, thin border) already used on
# the other sheets by copying its format onto the new header cells instead
# of re-deriving the same formatting from scratch.
$battingSheet.Range("A1:D1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

# Data row - force text storage (values look numeric / shouldn't become
# real numbers) then strip the temporary number format back off so the
# cells end up with no explicit style, matching the rest of the workbook.
$playerData = $playerInfo.Range("A2:D2")
$playerData.NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = "6922"
$playerInfo.Cells.Item(2,2).Value = "William George Jacks"
$playerInfo.Cells.Item(2,3).Value = "Right Handed"
$playerInfo.Cells.Item(2,4).Value = "Right Arm Off Break"
$playerData.ClearFormats()

# ---------------------------------------------------------------------
# 2. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE and replace the
#    full scorecard URLs with the bare match code.
# ---------------------------------------------------------------------
$battingSheet.Cells.Item(1,4).Value = "MATCH_CODE"

$battingCodes = $battingSheet.Range("D2:D3")
$battingCodes.NumberFormat = "@"
$battingSheet.Cells.Item(2,4).Value = "4711"
$battingSheet.Cells.Item(3,4).Value = "4713"
$battingCodes.ClearFormats()

# ---------------------------------------------------------------------
# 3. "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE rework.
# ---------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1,2).Value = "MATCH_CODE"

$bowlingCodes = $bowlingSheet.Range("B2:B3")
$bowlingCodes.NumberFormat = "@"
$bowlingSheet.Cells.Item(2,2).Value = "4711"
$bowlingSheet.Cells.Item(3,2).Value = "4713"
$bowlingCodes.ClearFormats()

Write-Output "Player Info sheet added; MATCH_CARD_LINK columns converted to MATCH_CODE."
